# Add data for 2025-07-18: update column L (year 2025 cumulative totals)
# for Citywide Totals, By Neighborhood, and affected neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3611  # was 3592
$ws.Range("L3").Value = 3782  # was 3765
$ws.Range("L4").Value = 940  # was 935
$ws.Range("L6").Value = 3299  # was 3280
$ws.Range("L7").Value = 11855  # was 11795

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 220  # was 218
$ws.Range("L3").Value = 256  # was 254
$ws.Range("L7").Value = 764  # was 760

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 89  # was 88
$ws.Range("L3").Value = 110  # was 107
$ws.Range("L7").Value = 274  # was 270

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 182  # was 181
$ws.Range("L6").Value = 183  # was 182
$ws.Range("L7").Value = 559  # was 557

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 138  # was 133
$ws.Range("L6").Value = 120  # was 119
$ws.Range("L7").Value = 426  # was 420

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 47  # was 46
$ws.Range("L7").Value = 396  # was 395
$ws.Range("L8").Value = 764  # was 760
$ws.Range("L10").Value = 74  # was 72
$ws.Range("L12").Value = 26  # was 25
$ws.Range("L20").Value = 298  # was 297
$ws.Range("L29").Value = 639  # was 637
$ws.Range("L31").Value = 115  # was 114
$ws.Range("L33").Value = 559  # was 557
$ws.Range("L37").Value = 426  # was 420
$ws.Range("L42").Value = 381  # was 379
$ws.Range("L43").Value = 88  # was 87
$ws.Range("L49").Value = 61  # was 59
$ws.Range("L50").Value = 56  # was 57
$ws.Range("L54").Value = 248  # was 245
$ws.Range("L55").Value = 112  # was 111
$ws.Range("L62").Value = 6  # was 5
$ws.Range("L63").Value = 43  # was 41
$ws.Range("L64").Value = 76  # was 75
$ws.Range("L66").Value = 30  # was 29
$ws.Range("L67").Value = 423  # was 420
$ws.Range("L73").Value = 101  # was 99
$ws.Range("L77").Value = 75  # was 72
$ws.Range("L78").Value = 149  # was 147
$ws.Range("L79").Value = 312  # was 308
$ws.Range("L83").Value = 274  # was 270
$ws.Range("L85").Value = 605  # was 603
$ws.Range("L89").Value = 167  # was 166
$ws.Range("L91").Value = 170  # was 169
$ws.Range("L92").Value = 35  # was 34
$ws.Range("L94").Value = 141  # was 140
$ws.Range("L96").Value = 120  # was 118
$ws.Range("L101").Value = 11855  # was 11795

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 36  # was 35
$ws.Range("L7").Value = 115  # was 114

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 122  # was 121
$ws.Range("L3").Value = 160  # was 158
$ws.Range("L7").Value = 423  # was 420

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L2").Value = 22  # was 21
$ws.Range("L6").Value = 26  # was 25
$ws.Range("L7").Value = 61  # was 59

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 53  # was 52
$ws.Range("L4").Value = 20  # was 19
$ws.Range("L6").Value = 116  # was 115
$ws.Range("L7").Value = 248  # was 245

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 194  # was 192
$ws.Range("L7").Value = 639  # was 637

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 116  # was 114
$ws.Range("L7").Value = 381  # was 379

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 32  # was 31
$ws.Range("L6").Value = 20  # was 19
$ws.Range("L7").Value = 74  # was 72

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 43  # was 42
$ws.Range("L6").Value = 46  # was 45
$ws.Range("L7").Value = 149  # was 147

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L4").Value = 7  # was 6
$ws.Range("L7").Value = 112  # was 111

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 33  # was 31
$ws.Range("L7").Value = 120  # was 118

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 72  # was 71
$ws.Range("L7").Value = 170  # was 169

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 113  # was 111
$ws.Range("L6").Value = 63  # was 61
$ws.Range("L7").Value = 312  # was 308

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 23  # was 22
$ws.Range("L7").Value = 76  # was 75

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 80  # was 79
$ws.Range("L7").Value = 298  # was 297

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 110  # was 109
$ws.Range("L7").Value = 396  # was 395

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 52  # was 51
$ws.Range("L7").Value = 141  # was 140

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L4").Value = 5  # was 6
$ws.Range("L7").Value = 56  # was 57

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L6").Value = 9  # was 8
$ws.Range("L7").Value = 30  # was 29

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 38  # was 37
$ws.Range("L3").Value = 29  # was 28
$ws.Range("L7").Value = 101  # was 99

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 15  # was 14
$ws.Range("L7").Value = 35  # was 34

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 42  # was 41
$ws.Range("L7").Value = 167  # was 166

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L4").Value = 3  # was 2
$ws.Range("L7").Value = 47  # was 46

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 25  # was 24
$ws.Range("L7").Value = 88  # was 87

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 176  # was 174
$ws.Range("L7").Value = 605  # was 603

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 24  # was 22
$ws.Range("L6").Value = 15  # was 14
$ws.Range("L7").Value = 75  # was 72

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L4").Value = 6  # was 5
$ws.Range("L7").Value = 26  # was 25

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("L2").Value = 4  # was 3
$ws.Range("L7").Value = 6  # was 5
